$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Types")
$ws.Range("C3:F3").ClearContents()
$ws.Range("C3:E3").Font.Bold = $true
